$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.245.87"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.246.43"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.62"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.19"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -5.64%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.611"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.22"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0947"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.16"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.07%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.42"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.850"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.233.70"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "42.140.54"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000100"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.34"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.57"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.03"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +35.12%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.51"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.61"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.65%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.61"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.69"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0822"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.33%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.77"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.43%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +9.36%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.73%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "13.75"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.67%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "62.53"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.08%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.86"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -6.23%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.72"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.16%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.102"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.78%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.43%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -8.61%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.10%  "
